$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All Fitness values in C2:C252 were overwritten with the constant 7573
$ws.Range("C2:C252").Value = 7573
